# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers in AD1:AF1, matching the style
# (bold, centered, top-aligned, thin border) used by the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerStyle = $ws.Range("A1")
$ws.Range("AD1:AF1").Font.Bold = $headerStyle.Font.Bold
$ws.Range("AD1:AF1").HorizontalAlignment = $headerStyle.HorizontalAlignment
$ws.Range("AD1:AF1").VerticalAlignment = $headerStyle.VerticalAlignment
$ws.Range("AD1:AF1").Borders.LineStyle = $headerStyle.Borders.LineStyle

# Data rows 2-43: the season record (69 wins, 93 losses, 0 ties) repeated
# for every player row.
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
